# Applies the gh-pages data-refresh commit (456a3b4):
# updates the "F" (want-to-go count) column across sheets 1, 2 and 4.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1889  # was 1886
$ws1.Range("F5").Value = 21  # was 20
$ws1.Range("F6").Value = 860  # was 854
$ws1.Range("F13").Value = 155  # was 154
$ws1.Range("F14").Value = 142  # was 141
$ws1.Range("F16").Value = 4451  # was 4444
$ws1.Range("F18").Value = 35  # was 32
$ws1.Range("F19").Value = 485  # was 483
$ws1.Range("F20").Value = 435  # was 433
$ws1.Range("F21").Value = 11  # was 9
$ws1.Range("F24").Value = 2121  # was 2085
$ws1.Range("F25").Value = 372  # was 371
$ws1.Range("F26").Value = 56  # was 52
$ws1.Range("F27").Value = 35  # was 33
$ws1.Range("F28").Value = 53  # was 52
$ws1.Range("F29").Value = 2147  # was 2134
$ws1.Range("F30").Value = 80  # was 79
$ws1.Range("F31").Value = 66  # was 64
$ws1.Range("F33").Value = 151  # was 150
$ws1.Range("F36").Value = 216  # was 215

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 35  # was 34

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1889  # was 1886
$ws4.Range("F5").Value = 21  # was 20
$ws4.Range("F6").Value = 860  # was 854
$ws4.Range("F13").Value = 155  # was 154
$ws4.Range("F14").Value = 142  # was 141
$ws4.Range("F16").Value = 35  # was 34
$ws4.Range("F17").Value = 4451  # was 4444
$ws4.Range("F19").Value = 35  # was 32
$ws4.Range("F20").Value = 485  # was 483
$ws4.Range("F21").Value = 435  # was 433
$ws4.Range("F22").Value = 11  # was 9
$ws4.Range("F25").Value = 2121  # was 2085
$ws4.Range("F26").Value = 372  # was 371
$ws4.Range("F27").Value = 56  # was 52
$ws4.Range("F28").Value = 35  # was 33
$ws4.Range("F29").Value = 53  # was 52
$ws4.Range("F30").Value = 2147  # was 2134
$ws4.Range("F31").Value = 80  # was 79
$ws4.Range("F32").Value = 66  # was 64
$ws4.Range("F34").Value = 151  # was 150
$ws4.Range("F37").Value = 216  # was 215

